$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RUN" column (A) for rows 7-12 was switched from "y" to "n"
# (skip these steps), matching rows 13-19 which already read "n".
$ws.Range("A7:A12").Value = "n"

# The saved selection moves to A13.
$ws.Range("A13").Select()
